$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price column values stay as text (avoid Excel auto-numeric coercion),
# matching the original inline-string cell type for D2:D51.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "29.201.75"
$ws.Range("E2").Value = "  -0.89%  "

$ws.Range("D3").Value = "1.860.47"
$ws.Range("E3").Value = "  -0.93%  "

$ws.Range("D4").Value = "1.000"
$ws.Range("E4").Value = "  -0.05%  "

$ws.Range("D5").Value = "0.7074"
$ws.Range("E5").Value = "  -0.83%  "

$ws.Range("D6").Value = "240.44"
$ws.Range("E6").Value = "  -0.71%  "

$ws.Range("E7").Value = "  +0.01%  "

$ws.Range("D8").Value = "0.3076"
$ws.Range("E8").Value = "  -1.26%  "

$ws.Range("D9").Value = "0.07653"
$ws.Range("E9").Value = "  -2.48%  "

$ws.Range("D10").Value = "24.77"
$ws.Range("E10").Value = "  -1.62%  "

$ws.Range("D11").Value = "0.08419"
$ws.Range("E11").Value = "  +1.88%  "

$ws.Range("D12").Value = "1.860.14"
$ws.Range("E12").Value = "  -0.17%  "

$ws.Range("D13").Value = "5.178"
$ws.Range("E13").Value = "  -1.85%  "

$ws.Range("D14").Value = "0.7094"
$ws.Range("E14").Value = "  -2.96%  "

$ws.Range("D15").Value = "91.02"
$ws.Range("E15").Value = "  -0.19%  "

$ws.Range("D16").Value = "29.209.97"
$ws.Range("E16").Value = "  -0.83%  "

$ws.Range("D17").Value = "5.932"
$ws.Range("E17").Value = "  +0.17%  "

$ws.Range("E18").Value = "  -1.78%  "

$ws.Range("D19").Value = "0.000007823"
$ws.Range("E19").Value = "  -0.70%  "

$ws.Range("D20").Value = "2.113.74"
$ws.Range("E20").Value = "  -0.10%  "

$ws.Range("E21").Value = "  -1.45%  "

$ws.Range("E22").Value = "  +0.08%  "

$ws.Range("D23").Value = "7.845"
$ws.Range("E23").Value = "  -1.38%  "

$ws.Range("E24").Value = "  -0.03%  "

$ws.Range("D25").Value = "0.1586"
$ws.Range("E25").Value = "  -0.22%  "

$ws.Range("D26").Value = "162.85"
$ws.Range("E26").Value = "  -0.56%  "

$ws.Range("D27").Value = "8.904"
$ws.Range("E27").Value = "  -1.25%  "

$ws.Range("D28").Value = "18.42"
$ws.Range("E28").Value = "  +0.61%  "

$ws.Range("E29").Value = "  +0.43%  "

$ws.Range("E30").Value = "  -3.57%  "

$ws.Range("D31").Value = "4.400"
$ws.Range("E31").Value = "  +0.61%  "

$ws.Range("D32").Value = "4.220"
$ws.Range("E32").Value = "  +2.15%  "

$ws.Range("D33").Value = "0.05128"
$ws.Range("E33").Value = "  -3.55%  "

$ws.Range("D34").Value = "0.8073"
$ws.Range("E34").Value = "  +11.53%  "

$ws.Range("D35").Value = "1.913"
$ws.Range("E35").Value = "  -1.15%  "

$ws.Range("D36").Value = "1.166"
$ws.Range("E36").Value = "  -2.98%  "

$ws.Range("D37").Value = "2.679"
$ws.Range("E37").Value = "  -0.06%  "

$ws.Range("D38").Value = "0.01845"
$ws.Range("E38").Value = "  -1.33%  "

$ws.Range("D39").Value = "2.691"
$ws.Range("E39").Value = "  -1.73%  "

$ws.Range("D40").Value = "1.177.12"
$ws.Range("E40").Value = "  -6.79%  "

$ws.Range("D41").Value = "6.184"
$ws.Range("E41").Value = "  +0.81%  "

$ws.Range("D42").Value = "0.8952"
$ws.Range("E42").Value = "  -1.72%  "

$ws.Range("D43").Value = "72.76"
$ws.Range("E43").Value = "  -1.76%  "

$ws.Range("E44").Value = "  -0.05%  "

$ws.Range("D45").Value = "101.83"
$ws.Range("E45").Value = "  -1.55%  "

$ws.Range("D46").Value = "2.008.49"
$ws.Range("E46").Value = "  +0.02%  "

$ws.Range("D47").Value = "0.5166"
$ws.Range("E47").Value = "  -3.07%  "

$ws.Range("D48").Value = "1.772"
$ws.Range("E48").Value = "  -0.05%  "

$ws.Range("E49").Value = "  -0.04%  "

$ws.Range("D50").Value = "9.259"
$ws.Range("E50").Value = "  +0.11%  "

$ws.Range("D51").Value = "1.002"
$ws.Range("E51").Value = "  +0.56%  "
